$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header + values (restoring "gk-aks-Digital" org exclude strings)
$ws.Range("D1").Value = "Exclude Strings"
$ws.Range("D2").Value = "gk-aks-Digital\nx-affected-action"
$ws.Range("D3").Value = "gk-aks-Digital\cf-tf-module-s3-website"
$ws.Range("D4").Value = "gk-aks-Digital\cf-tf-module-ecs"
$ws.Range("D5").Value = "gk-aks-Digital\testing123"
$ws.Range("D6").Value = "gk-aks-Digital\JamesRepo"
$ws.Range("D7").Value = "gk-aks-Digital\DeSilvaRepo"
$ws.Range("D8").Value = "gk-aks-Digital\RanjanRepo"
$ws.Range("D9").Value = "gk-aks-Digital\Testing12345"

# Match the new column width for D (Excel's ColumnWidth property is offset from
# the stored XML width by the default column padding; 32.1666... yields width="33")
$ws.Columns.Item(4).ColumnWidth = 32.16666666666666

# Update the selected cell to match the final state captured in the workbook
$ws.Range("G8").Select()
